$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new "id" column before column A: shift A:F -> B:G (values + the
# per-cell formatting that varies in this sheet), then fill column A with an
# "id" header and a 1..19 row counter.
# ---------------------------------------------------------------------------

# 1) Shift all values one column to the right (A1:F20 -> B1:G20). Using
#    Value2 on the whole block keeps shared-string reuse correct and does
#    NOT touch styles.xml at all.
$vals = $ws.Range("A1:F20").Value2
$ws.Range("B1:G20").Value2 = $vals

# 2) Move the per-row formatting that lived in column C (the revenue figures,
#    alternating style 2 / style 3) over to the new column D, without
#    creating any new style entries (PasteSpecial format-only reuses existing
#    cellXfs/font entries instead of allocating new ones).
$ws.Range("C2:C17").Copy() | Out-Null
$ws.Range("D2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# 3) Move the wrap-text style that lived on E9 over to the new F9.
$ws.Range("E9").Copy() | Out-Null
$ws.Range("F9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# 4) The old C2:C17 / E9 cells now hold shifted-in values (years / plain
#    text) that must go back to the default (unstyled) look. Reset them by
#    pasting the format of an always-default cell (A1) onto them.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C2:C17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("E9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.CutCopyMode = 0

# 5) Fill the new column A: header "id" then a 1-based row counter.
$ws.Range("A1").Value2 = "id"
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 1
}

# 6) Match the author's final selection.
$ws.Range("B7").Select() | Out-Null
